# BOI11.xlsx - "Add files via upload" re-sync
#
# The sheet tracks one row per (Account/Manager/Rep/Product) combination,
# grouped in blocks of four rows that share the same deal Status. This
# edit normalizes the Status column so every row in a block carries the
# same value, introduces a new "Pending" status (capitalized, distinct
# from the existing lowercase "pending") for the last block, and removes
# the trailing one-off row (row 18) that didn't belong to any block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize the Status column (H) so each 4-row deal block agrees ---

# Rows 2-5 (account block 1) -> "presented"
$excel.Union($ws.Range("H4"), $ws.Range("H5")).Value = "presented"

# Rows 6-9 (account block 2) -> "won"
$excel.Union($ws.Range("H7"), $ws.Range("H8"), $ws.Range("H9")).Value = "won"

# Rows 10-13 (account block 3) -> "declined"
$excel.Union($ws.Range("H11"), $ws.Range("H12"), $ws.Range("H13")).Value = "declined"

# Rows 14-17 (account block 4) -> "Pending" (new, capitalized status)
$excel.Union($ws.Range("H14"), $ws.Range("H15"), $ws.Range("H16"), $ws.Range("H17")).Value = "Pending"

# --- Drop the stray trailing row (account 729833 / Koepp Ltd / Monitor) ---
$ws.Rows(18).Delete()

# --- Leave the view parked on the tail of the table, like the source file ---
$ws.Range("H17").Select() | Out-Null
